# Daily update to data.
# Adds one new observation date (2020-03-26 11:00, serial 43916.458333333336)
# to the "longform" sheet (row 33) and its unpivoted equivalent rows
# (39 rows: 13 provinces/regions x 3 case types) to the "shortform" sheet
# (rows 366-404).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "longform" sheet - append row 33
# ---------------------------------------------------------------------------
$longform = $wb.Worksheets.Item("longform")

$dateSerial = 43916.458333333336

$longform.Cells.Item(33, 1).Value = "live"
$longform.Cells.Item(33, 2).Value = "govt_canada_ph"
$longform.Cells.Item(33, 3).Value = $dateSerial

# Copy the date-cell formatting from the row above so the new date cell
# keeps the existing "m/d/yyyy h:mm" style instead of creating a new one.
$longform.Range("C32").Copy()
$longform.Range("C33").PasteSpecial(-4122)

$longformValues = @(659,0,14,419,0,2,86,0,0,11,24,0,858,0,13,1339,0,6,67,0,0,26,0,0,68,0,0,5,0,0,13,0,0,3,0,0,1,0,0)
for ($i = 0; $i -lt $longformValues.Length; $i++) {
    $col = 4 + $i
    $longform.Cells.Item(33, $col).Value = $longformValues[$i]
}

# Update the selection to follow the newly entered data (matches Excel's
# behaviour of leaving the cursor on the last-edited cell) while keeping
# "shortform" as the tab that is actually active/selected in the workbook.
$null = $longform.Activate()
$null = $longform.Range("C33").Select()

# ---------------------------------------------------------------------------
# 2. "shortform" sheet - append rows 366-404 (unpivoted longform row 33)
# ---------------------------------------------------------------------------
$shortform = $wb.Worksheets.Item("shortform")
$null = $shortform.Activate()

$provinces = @("BC","BC","BC","AB","AB","AB","SK","SK","SK","MB","MB","MB","ON","ON","ON","QC","QC","QC","NL","NL","NL","NB","NB","NB","NS","NS","NS","PEI","PEI","PEI","Repat","Repat","Repat","YK","YK","YK","NT","NT","NT")
$caseTypes = @("conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths","conf","prob","deaths")
$cases     = @(659,0,14,419,0,2,86,0,0,11,24,0,858,0,13,1339,0,6,67,0,0,26,0,0,68,0,0,5,0,0,13,0,0,3,0,0,1,0,0)

$startRow = 366
for ($i = 0; $i -lt $provinces.Length; $i++) {
    $r = $startRow + $i
    $shortform.Cells.Item($r, 1).Value = "live"
    $shortform.Cells.Item($r, 2).Value = "govt_canada_ph"
    $shortform.Cells.Item($r, 3).Value = $dateSerial
    $shortform.Cells.Item($r, 4).Value = $provinces[$i]
    $shortform.Cells.Item($r, 5).Value = $caseTypes[$i]
    $shortform.Cells.Item($r, 6).Value = $cases[$i]
}

$lastRow = $startRow + $provinces.Length - 1
$null = $shortform.Range("C366:C" + $lastRow).Select()
